$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: percentage split of the copy/paste totals (mirrors row 4's pattern)
$ws.Range("A10").Formula = "=A9*100/D8"
$ws.Range("B10").Formula = "=B9*100/D8"
$ws.Range("A10:B10").Style = $ws.Range("A2").Style

# Row 11 left blank but formatted the same as the block above
$ws.Range("A11:B11").Style = $ws.Range("A2").Style

# Row 12 - label marking this block as the pasted copy
$ws.Range("A12").Style = $ws.Range("A2").Style
$ws.Range("B12").Value = "copy^"
$ws.Range("B12").Style = $ws.Range("A2").Style

# Update selection like Excel would leave it after this editing session
$ws.Range("E9").Select()

$wb.Save()
